# Update "想去人数" (column F) values across the four worksheets to match
# the latest generated snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) pairs for column F.
$updates = @{
    "展览" = @{
        5  = 2982
        10 = 7023
        11 = 46
        12 = 119
        13 = 382
        14 = 618
        15 = 1523
        16 = 1138
        17 = 2274
        18 = 1525
        19 = 135
        21 = 144
        22 = 7
        24 = 358
        25 = 38
        26 = 38
        27 = 1777
        31 = 1684
        32 = 1246
        35 = 10
        42 = 33
        47 = 331
    }
    "演出" = @{
        6  = 184
        15 = 61
        20 = 66
        37 = 41
    }
    "本地生活" = @{
        6  = 1714
        8  = 2770
        9  = 1046
        10 = 965
        12 = 300
        13 = 1567
        14 = 7459
    }
    "全部类型" = @{
        5  = 2982
        6  = 1714
        7  = 2770
        8  = 7023
        9  = 1046
        10 = 46
        11 = 382
        13 = 1567
        14 = 618
        15 = 1523
        16 = 1138
        17 = 2274
        18 = 1525
        19 = 135
        21 = 144
        22 = 7
        24 = 38
        25 = 38
        26 = 1777
        29 = 1684
        30 = 1246
        35 = 66
        42 = 33
        50 = 41
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
